$re = New-Object System.Text.RegularExpressions.Regex
Write-Output ("re type: " + $re.GetType().FullName)
